# Update balance/cg values on "GLOBAL RESULTS" sheet following updated
# climb performance / operating conditions computations.
$wb = $excel.ActiveWorkbook

$global = $wb.Worksheets.Item("GLOBAL RESULTS")

$global.Range("C2").Value = 5.358157033106666
$global.Range("C3").Value = 12.521056690418195
$global.Range("C4").Value = 0.7221121739534344

$global.Range("C6").Value = 5.103819398632732
$global.Range("C7").Value = 11.935465713143788
$global.Range("C8").Value = 0.784349506445851

$global.Range("C10").Value = 5.103819398632732
$global.Range("C11").Value = 11.935465713143788
$global.Range("C12").Value = 0.784349506445851

$global.Range("C14").Value = 5.125710773137314
$global.Range("C15").Value = 11.985868758750996
$global.Range("C16").Value = 0.4871321042486384

$global.Range("C18").Value = 5.099577445159801
$global.Range("C19").Value = 11.925698972626193
$global.Range("C20").Value = 0.7143954860380879

# Update landing gears Xcg BRF value
$landingGears = $wb.Worksheets.Item("LANDING GEARS")
$landingGears.Range("C2").Value = 12.308548373872053

$wb.Save()
